# Apply the change: insert a new row before existing row 4 on the first
# worksheet (strategy_id-0) and populate it with the new
# "climate_change_factor_gnrl_hydropower_availability" variable, shifting
# all the following rows (elasticity_gnrl_rate_occupancy_to_gdppc,
# frac_gnrl_eating_red_meat, limit_gnrl_annual_emissions_mt_ch4/co2/n2o,
# occrateinit_gnrl_occupancy, population_gnrl_rural/urban) down by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 4; this shifts rows 4-11 down to 5-12 and
# automatically grows the sheet dimension from A1:AS11 to A1:AS12.
$ws.Rows.Item(4).Insert()

# Populate the newly inserted row 4 with the new variable definition.
# (Columns C:G are metadata flags that are blank for every row in this
# sheet, so row 4 simply keeps them blank like its neighbours.)
$ws.Range("A4").Value = "General"
$ws.Range("B4").Value = "climate_change_factor_gnrl_hydropower_availability"
$ws.Range("H4").Value = 1
$ws.Range("I4").Value = 0.5

# Columns J4:AS4 all carry the value 1 for this new row.
$ws.Range("J4:AS4").Value = 1
